{"js": "// The document has a stray one-character paragraph (\"j\") left over right\n// after the \"Answer: java\" line for question 1 (a leftover keystroke that\n// was never meant to be there). The fix removes that leftover run of text\n// so the paragraph goes back to being empty, the same as the blank\n// paragraphs surrounding it.\n\n// Search the body for the exact, whole-word \"j\" (case sensitive) \u2014 this\n// uniquely matches the stray run and nothing else (the document elsewhere\n// only ever has \"j\" as part of longer words like \"Java\", \"JavaScript\",\n// \"JSX\", etc., which matchWholeWord excludes).\nconst results = context.document.body.search(\"j\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nresults.load(\"items\");\nawait context.sync();\n\n// Delete the matched run of text in place; because it is the only content\n// in its paragraph, the paragraph itself remains (now empty) exactly like\n// the diff shows \u2014 only the run carrying \"j\" disappears.\nif (results.items.length > 0) {\n  results.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# The document has a stray one-character paragraph (\"j\") left over right\n# after the \"Answer: java\" line for question 1 (a leftover keystroke that\n# was never meant to be there). The fix removes that leftover run of text\n# so the paragraph goes back to being empty, the same as the blank\n# paragraphs surrounding it.\n\n$d = $word.ActiveDocument\n\n# Search the whole document body for the exact, whole-word \"j\"\n# (case sensitive) -- this uniquely matches the stray run and nothing else\n# (the document elsewhere only ever has \"j\" as part of longer words like\n# \"Java\", \"JavaScript\", \"JSX\", etc., which MatchWholeWord excludes).\n$range = $d.Content\n$range.Find.ClearFormatting()\n$found = $range.Find.Execute(\"j\", $true, $true)\n\n# Delete the matched run of text in place; because it is the only content\n# in its paragraph, the paragraph itself remains (now empty) exactly like\n# the diff shows -- only the run carrying \"j\" disappears.\nif ($found) {\n    $range.Delete()\n}\n"}
